# Advances in ktk.cycles ; Started a tutorial for TimeSeries.
#
# Updates the pyKTK Matlab->Python conversion tracking sheet:
#   - kinematics.readc3dfile note: switched from the forked py-c3d lib to ezc3d
#   - TimeSeries.findrepeatablecycles note: now planned for ktk.cycles
#   - TimeSeries.resample: no longer TODO, Python name is simply "resample"
#   - TimeSeries.sortevents: no longer TODO, not required (sorted(ts.events))
#   - moves the sheet's scroll/selection down to around row 45 / cell F68

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 31: kinematics / readc3dfile ---------------------------------
$ws.Range("F31").Value = "Now using ezc3d to be up to date with different c3d formats"

# --- Row 50: TimeSeries / findrepeatablecycles -------------------------
$ws.Range("F50").Value = "Will be in ktk.cycles"

# --- Row 65: TimeSeries / resample --------------------------------------
$ws.Range("D65").ClearContents()
$ws.Range("E65").ClearContents()
$ws.Range("C65").Value = "resample"

# --- Row 66: TimeSeries / sortevents -------------------------------------
$ws.Range("D66").ClearContents()
$ws.Range("E66").ClearContents()
$ws.Range("C66").Value = "N/A"
$ws.Range("F66").Value = "Not required, we can use sorted(ts.events) to get the events sorted."

# --- Scroll / selection state -------------------------------------------
$ws.Range("F68").Select()
